# GSC_MIxS_6.form.xlsx — MimsSoil sheet: reorder the header row (row 1,
# columns A:AR) and move the three affected data-validation dropdowns so
# they keep tracking the columns whose headers they were attached to.
# (Commit message: "generates GH pages! removes temporary elements!")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MimsSoil")

# --- 1. Rewrite the header row with the new column order -------------------
$ws.Range("A1").Value  = "collection_date"
$ws.Range("B1").Value  = "depth"
$ws.Range("C1").Value  = "tax_class"
$ws.Range("D1").Value  = "geo_loc_name"
$ws.Range("E1").Value  = "project_name"
$ws.Range("F1").Value  = "nucl_acid_ext"
$ws.Range("G1").Value  = "lib_reads_seqd"
$ws.Range("H1").Value  = "lat_lon"
$ws.Range("I1").Value  = "env_local_scale"
$ws.Range("J1").Value  = "samp_name"
$ws.Range("K1").Value  = "elev"
$ws.Range("L1").Value  = "sim_search_meth"
$ws.Range("M1").Value  = "temp"
$ws.Range("N1").Value  = "samp_taxon_id"
$ws.Range("O1").Value  = "samp_mat_process"
$ws.Range("P1").Value  = "lib_screen"
$ws.Range("Q1").Value  = "seq_meth"
$ws.Range("R1").Value  = "samp_size"
$ws.Range("S1").Value  = "source_mat_id"
$ws.Range("T1").Value  = "mid"
$ws.Range("U1").Value  = "assembly_qual"
$ws.Range("V1").Value  = "size_frac"
$ws.Range("W1").Value  = "env_medium"
$ws.Range("X1").Value  = "samp_collect_device"
$ws.Range("Y1").Value  = "feat_pred"
$ws.Range("Z1").Value  = "lib_size"
$ws.Range("AA1").Value = "env_broad_scale"
$ws.Range("AB1").Value = "lib_vector"
$ws.Range("AC1").Value = "assembly_name"
$ws.Range("AD1").Value = "samp_vol_we_dna_ext"
$ws.Range("AE1").Value = "adapters"
$ws.Range("AF1").Value = "number_contig"
$ws.Range("AG1").Value = "neg_cont_type"
$ws.Range("AH1").Value = "nucl_acid_amp"
$ws.Range("AI1").Value = "alt"
$ws.Range("AJ1").Value = "lib_layout"
$ws.Range("AK1").Value = "annot"
$ws.Range("AL1").Value = "experimental_factor"
$ws.Range("AM1").Value = "pos_cont_type"
$ws.Range("AN1").Value = "ref_biomaterial"
$ws.Range("AO1").Value = "assembly_software"
$ws.Range("AP1").Value = "rel_to_oxygen"
$ws.Range("AQ1").Value = "ref_db"
$ws.Range("AR1").Value = "samp_collect_method"
# columns AS1 onward (associated_resource, sop, ...) are unchanged.

# --- 2. Re-home the 3 data validations whose target column moved -----------
# neg_cont_type: was column N -> now column AG
# lib_layout:    was column Y -> now column AJ
# rel_to_oxygen: was column AH -> now column AP
# The other 5 validations (tillage/soil_horizon/fao_class/profile_position/
# drainage_class) keep their columns. All 8 are deleted and re-added, in the
# final desired order, so the moved ones land back in their original slot
# (positions 1-3) instead of being appended at the end.
$allRanges = @(
    "N2:N1048576","Y2:Y1048576","AH2:AH1048576",
    "BB2:BB1048576","BF2:BF1048576","BS2:BS1048576","CA2:CA1048576","CB2:CB1048576"
)
foreach ($r in $allRanges) {
    $ws.Range($r).Validation.Delete()
}

function Add-ListValidation($range, $formula) {
    $v = $range.Validation
    $v.Add(3, 1, 1, $formula)
    $v.IgnoreBlank = $true
    $v.InCellDropdown = $true
    $v.ShowInput = $false
    $v.ShowError = $false
}

Add-ListValidation $ws.Range("AG2:AG1048576") """DNA-free PCR mix,distilled water,empty collection device,empty collection tube,phosphate buffer,sterile swab,sterile syringe"""
Add-ListValidation $ws.Range("AJ2:AJ1048576") """other,paired,single,vector"""
Add-ListValidation $ws.Range("AP2:AP1048576") """aerobe,anaerobe,facultative,microaerophilic,microanaerobe,obligate aerobe,obligate anaerobe"""
Add-ListValidation $ws.Range("BB2:BB1048576") """chisel,cutting disc,disc plough,drill,mouldboard,ridge till,strip tillage,tined,zonal tillage"""
Add-ListValidation $ws.Range("BF2:BF1048576") """A horizon,B horizon,C horizon,E horizon,O horizon,Permafrost,R layer"""
Add-ListValidation $ws.Range("BS2:BS1048576") """Acrisols,Andosols,Arenosols,Cambisols,Chernozems,Ferralsols,Fluvisols,Gleysols,Greyzems,Gypsisols,Histosols,Kastanozems,Lithosols,Luvisols,Nitosols,Phaeozems,Planosols,Podzols,Podzoluvisols,Rankers,Regosols,Rendzinas,Solonchaks,Solonetz,Vertisols,Yermosols"""
Add-ListValidation $ws.Range("CA2:CA1048576") """backslope,footslope,shoulder,summit,toeslope"""
Add-ListValidation $ws.Range("CB2:CB1048576") """excessively drained,moderately well,poorly,somewhat poorly,very poorly,well"""

Write-Host "MimsSoil header reorder + validation remap complete"
